# Append the Friday, Jan 13 departures appended to the "Main Data" table
# (rows 115-129), matching the new rows added to DataBase/POZ_Departures.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFlights = @(
    ,(114, "Friday, Jan 13", "5:40 AM", "LO3942", "Warsaw", "(WAW)", "LOT ", "E75S", "(SP-LIC)", "5:50 AM", $null, "0 hours, 10 minutes", $null)
    ,(115, "Friday, Jan 13", "5:50 AM", "FR7938", "Edinburgh", "(EDI)", "Ryanair ", "B738", "(SP-RSM)", "5:57 AM", $null, "0 hours, 7 minutes", $null)
    ,(116, "Friday, Jan 13", "6:30 AM", "LH1381", "Frankfurt", "(FRA)", "Lufthansa ", "CRJ9", "(D-ACNN)", "6:32 AM", $null, "0 hours, 2 minutes", $null)
    ,(117, "Friday, Jan 13", "6:45 AM", "FR7942", "Manchester", "(MAN)", "Ryanair ", "B738", "(SP-RSX)", "6:44 AM", $null, "0 hours, -1 minutes", $null)
    ,(118, "Friday, Jan 13", "7:35 AM", "FR1751", "London", "(STN)", "Ryanair ", "B738", "(SP-RKR)", "7:36 AM", $null, "0 hours, 1 minutes", $null)
    ,(119, "Friday, Jan 13", "8:05 AM", "E47013", "Dubai", "(DWC)", "Enter Air ", "B738", "(SP-ENL)", "8:07 AM", $null, "0 hours, 2 minutes", $null)
    ,(120, "Friday, Jan 13", "8:05 AM", "LH1641", "Munich", "(MUC)", "Lufthansa ", "CRJ9", "(D-ACKF)", "8:09 AM", $null, "0 hours, 4 minutes", $null)
    ,(121, "Friday, Jan 13", "10:45 AM", "LH1389", "Frankfurt", "(FRA)", "Lufthansa ", "CRJ9", "(D-ACNO)", "10:56 AM", $null, "0 hours, 11 minutes", $null)
    ,(122, "Friday, Jan 13", "11:00 AM", "LPR42", "Warsaw", "(WAW)", "Polish Medical Air Rescue ", "LJ75", "(SP-MXS)", "11:01 AM", $null, "0 hours, 1 minutes", $null)
    ,(123, "Friday, Jan 13", "11:55 AM", "FR7891", "Venice", "(TSF)", "Ryanair ", "B738", "(SP-RSM)", "11:55 AM", $null, "0 hours, 0 minutes", $null)
    ,(124, "Friday, Jan 13", "12:05 PM", "FR7944", "Leeds", "(LBA)", "Ryanair ", "B738", "(SP-RSX)", "12:04 PM", $null, "0 hours, -1 minutes", $null)
    ,(125, "Friday, Jan 13", "12:15 PM", "W61925", "Eindhoven", "(EIN)", "Wizz Air ", "A320", "(HA-LYE)", "12:21 PM", $null, "0 hours, 6 minutes", $null)
    ,(126, "Friday, Jan 13", "12:40 PM", "FR9315", "Helsinki", "(HEL)", "Ryanair ", "B738", "(SP-RKR)", "12:39 PM", $null, "0 hours, -1 minutes", $null)
    ,(127, "Friday, Jan 13", "12:55 PM", "LH1637", "Munich", "(MUC)", "Lufthansa ", "CRJ9", "(D-ACNM)", "12:57 PM", $null, "0 hours, 2 minutes", $null)
    ,(128, "Friday, Jan 13", "1:10 PM", "W61921", "Paris", "(BVA)", "Wizz Air ", "A320", "(HA-LWV)", "1:15 PM", $null, "0 hours, 5 minutes", $null)
)

$startRow = 115
for ($i = 0; $i -lt $newFlights.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newFlights[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $val = $rowData[$col - 1]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $col).Value = $val
        }
    }
}
